# Regenerate the "K" column (column G) values for the save_data sheet.
# The author's commit switched the statistic stored in column G from a
# "Strike#" derived figure to a recomputed "K" value (std/mean based calc),
# and rewrote the stored s_vals. The net effect on this worksheet is a
# per-row update of column G (rows 2-58). Row 12 and row 51 already held
# the correct value (0) and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 3
    4  = 1
    5  = 0
    6  = 2
    7  = 0
    8  = 1
    9  = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 3
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 1
    25 = 3
    26 = 0
    27 = 1
    28 = 1
    29 = 2
    30 = 1
    31 = 0
    32 = 1
    33 = 0
    34 = 0
    35 = 1
    36 = 1
    37 = 0
    38 = 1
    39 = 2
    40 = 1
    41 = 1
    42 = 0
    43 = 1
    44 = 1
    45 = 0
    46 = 1
    47 = 3
    48 = 0
    49 = 1
    50 = 1
    51 = 0
    52 = 3
    53 = 1
    54 = 2
    55 = 2
    56 = 1
    57 = 0
    58 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
